# Apply the GBIF species-name corrections described in the commit:
# "always keep family infos, keep orig names also if not changed by GBIF"
#
# Column A (Species) text is updated for specific rows; Column B (Family)
# is left untouched throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3   = "Acer campestre"
    4   = "Acer negundo"
    5   = "Acer platanoides"
    6   = "Acer pseudoplatanus"
    7   = "Acer species"
    11  = "Allium species"
    17  = "Apiaceae species"
    22  = "Asteraceae species"
    27  = "Betula species"
    32  = "Brassicaceae species"
    66  = "Draba species"
    77  = "Festuca ovina"
    78  = "Festuca rubra"
    79  = "Festuca rupicola"
    80  = "Festuca species"
    91  = "Geranium molle"
    92  = "Geranium pratense"
    93  = "Geranium pusillum"
    94  = "Geranium pyrenaicum"
    95  = "Geranium rotundifolium"
    96  = "Geranium species"
    123 = "Medicago falcata"
    124 = "Medicago lupulina"
    125 = "Medicago species"
    145 = "Poaceae species"
    151 = "Prunus avium"
    152 = "Prunus mahaleb"
    153 = "Prunus species"
    159 = "Rubus caesius"
    160 = "Rubus idaeus"
    161 = "Rubus species"
    167 = "Senecio jacobaea"
    168 = "Senecio species"
    202 = "Triticum species"
    213 = "Vicia cracca"
    214 = "Vicia hirsuta"
    215 = "Vicia sativa"
    216 = "Vicia sepium"
    217 = "Vicia species"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}
